$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2169893333333333
$ws.Range("H2").Value = 0.650968
$ws.Range("I2").Value = 0.2138091362408864
$ws.Range("J2").Value = 0.2138091362408864
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 2.044875435241778
$ws.Range("R2").Value = 18.403878917176
$ws.Range("S2").Value = 0.01422757894466793
$ws.Range("T2").Value = 0.01422757894466793
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2169893333333333
$ws.Range("H3").Value = 0.650968
$ws.Range("I3").Value = 0.2138091362408864
$ws.Range("J3").Value = 0.2138091362408864
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 10.97804687864356
$ws.Range("R3").Value = 98.802421907792
$ws.Range("S3").Value = 0.07638168366265258
$ws.Range("T3").Value = 0.07638168366265258
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2169893333333333
$ws.Range("H4").Value = 0.650968
$ws.Range("I4").Value = 0.2138091362408864
$ws.Range("J4").Value = 0.2138091362408864
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 5.824159197198222
$ws.Range("R4").Value = 52.417432774784
$ws.Range("S4").Value = 0.04052260755660845
$ws.Range("T4").Value = 0.04052260755660846
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2169893333333333
$ws.Range("H5").Value = 0.650968
$ws.Range("I5").Value = 0.2138091362408864
$ws.Range("J5").Value = 0.2138091362408864
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 11.88288682925067
$ws.Range("R5").Value = 106.945981463256
$ws.Range("S5").Value = 0.08267726607695741
$ws.Range("T5").Value = 0.08267726607695741
$ws.Range("I6").Value = 0.3272432505578689
$ws.Range("J6").Value = 0.3272432505578689
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 3.129761880991556
$ws.Range("R6").Value = 28.167856928924
$ws.Range("S6").Value = 0.02177586637914442
$ws.Range("T6").Value = 0.02177586637914442
$ws.Range("I7").Value = 0.3272432505578689
$ws.Range("J7").Value = 0.3272432505578689
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("S7").Value = 0.1169051560859796
$ws.Range("T7").Value = 0.1169051560859796
$ws.Range("I8").Value = 0.3272432505578689
$ws.Range("J8").Value = 0.3272432505578689
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 8.914103583068444
$ws.Range("R8").Value = 80.22693224761601
$ws.Range("S8").Value = 0.06202143674050155
$ws.Range("T8").Value = 0.06202143674050156
$ws.Range("I9").Value = 0.3272432505578689
$ws.Range("J9").Value = 0.3272432505578689
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 18.18722333564933
$ws.Range("R9").Value = 163.685010020844
$ws.Range("S9").Value = 0.1265407913522433
$ws.Range("T9").Value = 0.1265407913522433
$ws.Range("G10").Value = 0.4136746666666666
$ws.Range("H10").Value = 1.241024
$ws.Range("I10").Value = 0.4076118480389355
$ws.Range("J10").Value = 0.4076118480389355
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 3.898408972707556
$ws.Range("R10").Value = 35.085680754368
$ws.Range("S10").Value = 0.02712386312726212
$ws.Range("T10").Value = 0.02712386312726212
$ws.Range("G11").Value = 0.4136746666666666
$ws.Range("H11").Value = 1.241024
$ws.Range("I11").Value = 0.4076118480389355
$ws.Range("J11").Value = 0.4076118480389355
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 20.92886232429511
$ws.Range("R11").Value = 188.359760918656
$ws.Range("S11").Value = 0.1456162247387886
$ws.Range("T11").Value = 0.1456162247387886
$ws.Range("G12").Value = 0.4136746666666666
$ws.Range("H12").Value = 1.241024
$ws.Range("I12").Value = 0.4076118480389355
$ws.Range("J12").Value = 0.4076118480389355
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 11.10334354921244
$ws.Range("R12").Value = 99.930091942912
$ws.Range("S12").Value = 0.07725345719041864
$ws.Range("T12").Value = 0.07725345719041866
$ws.Range("G13").Value = 0.4136746666666666
$ws.Range("H13").Value = 1.241024
$ws.Range("I13").Value = 0.4076118480389355
$ws.Range("J13").Value = 0.4076118480389355
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 22.65387506664533
$ws.Range("R13").Value = 203.884875599808
$ws.Range("S13").Value = 0.1576183029824661
$ws.Range("T13").Value = 0.1576183029824661
$ws.Range("G14").Value = 0.05209933333333334
$ws.Range("H14").Value = 0.156298
$ws.Range("I14").Value = 0.05133576516230915
$ws.Range("J14").Value = 0.05133576516230916
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 0.4909764239984445
$ws.Range("R14").Value = 4.418787815986001
$ws.Range("S14").Value = 0.003416054451053981
$ws.Range("T14").Value = 0.003416054451053981
$ws.Range("G15").Value = 0.05209933333333334
$ws.Range("H15").Value = 0.156298
$ws.Range("I15").Value = 0.05133576516230915
$ws.Range("J15").Value = 0.05133576516230916
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 2.635838890756889
$ws.Range("R15").Value = 23.722550016812
$ws.Range("S15").Value = 0.01833931067749148
$ws.Range("T15").Value = 0.01833931067749148
$ws.Range("G16").Value = 0.05209933333333334
$ws.Range("H16").Value = 0.156298
$ws.Range("I16").Value = 0.05133576516230915
$ws.Range("J16").Value = 0.05133576516230916
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 1.398385841091556
$ws.Range("R16").Value = 12.585472569824
$ws.Range("S16").Value = 0.009729514378406908
$ws.Range("T16").Value = 0.009729514378406912
$ws.Range("G17").Value = 0.05209933333333334
$ws.Range("H17").Value = 0.156298
$ws.Range("I17").Value = 0.05133576516230915
$ws.Range("J17").Value = 0.05133576516230916
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 2.853091773540667
$ws.Range("R17").Value = 25.677825961866
$ws.Range("S17").Value = 0.01985088565535678
$ws.Range("T17").Value = 0.01985088565535678
